$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title placeholder: "Robot DESIGN Lesson" -> "Technic Basics"
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Technic Basics"

# Subtitle placeholder: "Basic Building Techniques for LEGO Robots" -> "Seshan Brothers"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Seshan Brothers"
